$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).NumberFormat = "@"
$ws.Cells.Item(2,4).Value2 = '29.046.40'
$ws.Cells.Item(2,4).Style = "Normal"
$ws.Cells.Item(2,5).Value2 = '  -2.04%  '

$ws.Cells.Item(3,4).NumberFormat = "@"
$ws.Cells.Item(3,4).Value2 = '1.992.25'
$ws.Cells.Item(3,4).Style = "Normal"
$ws.Cells.Item(3,5).Value2 = '  -1.08%  '

$ws.Cells.Item(4,4).NumberFormat = "@"
$ws.Cells.Item(4,4).Value2 = '1.015'
$ws.Cells.Item(4,4).Style = "Normal"
$ws.Cells.Item(4,5).Value2 = '  +0.39%  '

$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value2 = '330.22'
$ws.Cells.Item(5,4).Style = "Normal"
$ws.Cells.Item(5,5).Value2 = '  -0.68%  '

$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value2 = '1.012'
$ws.Cells.Item(6,4).Style = "Normal"
$ws.Cells.Item(6,5).Value2 = '  +0.34%  '

$ws.Cells.Item(7,4).NumberFormat = "@"
$ws.Cells.Item(7,4).Value2 = '0.4959'
$ws.Cells.Item(7,4).Style = "Normal"
$ws.Cells.Item(7,5).Value2 = '  -1.97%  '

$ws.Cells.Item(8,4).NumberFormat = "@"
$ws.Cells.Item(8,4).Value2 = '0.4182'
$ws.Cells.Item(8,4).Style = "Normal"
$ws.Cells.Item(8,5).Value2 = '  -2.10%  '

$ws.Cells.Item(9,4).NumberFormat = "@"
$ws.Cells.Item(9,4).Value2 = '55.19'
$ws.Cells.Item(9,4).Style = "Normal"
$ws.Cells.Item(9,5).Value2 = '  +1.82%  '

$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value2 = '0.08873'
$ws.Cells.Item(10,4).Style = "Normal"
$ws.Cells.Item(10,5).Value2 = '  -3.95%  '

$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(11,4).Value2 = '1.088'
$ws.Cells.Item(11,4).Style = "Normal"
$ws.Cells.Item(11,5).Value2 = '  -3.66%  '

$ws.Cells.Item(12,4).NumberFormat = "@"
$ws.Cells.Item(12,4).Value2 = '22.84'
$ws.Cells.Item(12,4).Style = "Normal"
$ws.Cells.Item(12,5).Value2 = '  -3.45%  '

$ws.Cells.Item(13,4).NumberFormat = "@"
$ws.Cells.Item(13,4).Value2 = '2.004.62'
$ws.Cells.Item(13,4).Style = "Normal"
$ws.Cells.Item(13,5).Value2 = '  +0.62%  '

$ws.Cells.Item(14,4).NumberFormat = "@"
$ws.Cells.Item(14,4).Value2 = '7.985'
$ws.Cells.Item(14,4).Style = "Normal"
$ws.Cells.Item(14,5).Value2 = '  -2.12%  '

$ws.Cells.Item(15,4).NumberFormat = "@"
$ws.Cells.Item(15,4).Value2 = '6.400'
$ws.Cells.Item(15,4).Style = "Normal"
$ws.Cells.Item(15,5).Value2 = '  -2.45%  '

$ws.Cells.Item(16,5).Value2 = '  +0.42%  '

$ws.Cells.Item(17,4).NumberFormat = "@"
$ws.Cells.Item(17,4).Value2 = '92.13'
$ws.Cells.Item(17,4).Style = "Normal"
$ws.Cells.Item(17,5).Value2 = '  -4.04%  '

$ws.Cells.Item(18,4).NumberFormat = "@"
$ws.Cells.Item(18,4).Value2 = '0.00001105'
$ws.Cells.Item(18,4).Style = "Normal"
$ws.Cells.Item(18,5).Value2 = '  -2.28%  '

$ws.Cells.Item(19,4).NumberFormat = "@"
$ws.Cells.Item(19,4).Value2 = '0.06718'
$ws.Cells.Item(19,4).Style = "Normal"
$ws.Cells.Item(19,5).Value2 = '  +0.91%  '

$ws.Cells.Item(20,4).NumberFormat = "@"
$ws.Cells.Item(20,4).Value2 = '19.42'
$ws.Cells.Item(20,4).Style = "Normal"
$ws.Cells.Item(20,5).Value2 = '  -2.93%  '

$ws.Cells.Item(21,4).NumberFormat = "@"
$ws.Cells.Item(21,4).Value2 = '1.012'
$ws.Cells.Item(21,4).Style = "Normal"
$ws.Cells.Item(21,5).Value2 = '  +0.47%  '

$ws.Cells.Item(22,4).NumberFormat = "@"
$ws.Cells.Item(22,4).Value2 = '5.966'
$ws.Cells.Item(22,4).Style = "Normal"
$ws.Cells.Item(22,5).Value2 = '  -0.74%  '

$ws.Cells.Item(23,4).NumberFormat = "@"
$ws.Cells.Item(23,4).Value2 = '29.102.57'
$ws.Cells.Item(23,4).Style = "Normal"
$ws.Cells.Item(23,5).Value2 = '  -2.00%  '

$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(24,4).Value2 = '11.94'
$ws.Cells.Item(24,4).Style = "Normal"
$ws.Cells.Item(24,5).Value2 = '  -1.05%  '

$ws.Cells.Item(25,4).NumberFormat = "@"
$ws.Cells.Item(25,4).Value2 = '2.322'
$ws.Cells.Item(25,4).Style = "Normal"
$ws.Cells.Item(25,5).Value2 = '  +1.88%  '

$ws.Cells.Item(26,4).NumberFormat = "@"
$ws.Cells.Item(26,4).Value2 = '2.244.16'
$ws.Cells.Item(26,4).Style = "Normal"
$ws.Cells.Item(26,5).Value2 = '  +0.49%  '

$ws.Cells.Item(27,4).NumberFormat = "@"
$ws.Cells.Item(27,4).Value2 = '20.79'
$ws.Cells.Item(27,4).Style = "Normal"
$ws.Cells.Item(27,5).Value2 = '  -0.41%  '

$ws.Cells.Item(28,4).NumberFormat = "@"
$ws.Cells.Item(28,4).Value2 = '156.70'
$ws.Cells.Item(28,4).Style = "Normal"
$ws.Cells.Item(28,5).Value2 = '  -1.83%  '

$ws.Cells.Item(29,4).NumberFormat = "@"
$ws.Cells.Item(29,4).Value2 = '6.259'
$ws.Cells.Item(29,4).Style = "Normal"
$ws.Cells.Item(29,5).Value2 = '  -4.21%  '

$ws.Cells.Item(30,4).NumberFormat = "@"
$ws.Cells.Item(30,4).Value2 = '2.240'
$ws.Cells.Item(30,4).Style = "Normal"
$ws.Cells.Item(30,5).Value2 = '  -4.92%  '

$ws.Cells.Item(31,4).NumberFormat = "@"
$ws.Cells.Item(31,4).Value2 = '126.76'
$ws.Cells.Item(31,4).Style = "Normal"
$ws.Cells.Item(31,5).Value2 = '  -1.87%  '

$ws.Cells.Item(32,4).NumberFormat = "@"
$ws.Cells.Item(32,4).Value2 = '1.039'
$ws.Cells.Item(32,4).Style = "Normal"
$ws.Cells.Item(32,5).Value2 = '  -2.45%  '

$ws.Cells.Item(33,4).NumberFormat = "@"
$ws.Cells.Item(33,4).Value2 = '0.09865'
$ws.Cells.Item(33,4).Style = "Normal"
$ws.Cells.Item(33,5).Value2 = '  -1.45%  '

$ws.Cells.Item(34,4).NumberFormat = "@"
$ws.Cells.Item(34,4).Value2 = '1.524'
$ws.Cells.Item(34,4).Style = "Normal"
$ws.Cells.Item(34,5).Value2 = '  -5.12%  '

$ws.Cells.Item(35,4).NumberFormat = "@"
$ws.Cells.Item(35,4).Value2 = '5.822'
$ws.Cells.Item(35,4).Style = "Normal"
$ws.Cells.Item(35,5).Value2 = '  -1.35%  '

$ws.Cells.Item(36,4).NumberFormat = "@"
$ws.Cells.Item(36,4).Value2 = '3.760'
$ws.Cells.Item(36,4).Style = "Normal"
$ws.Cells.Item(36,5).Value2 = '  -1.24%  '

$ws.Cells.Item(37,4).NumberFormat = "@"
$ws.Cells.Item(37,4).Value2 = '0.02410'
$ws.Cells.Item(37,4).Style = "Normal"
$ws.Cells.Item(37,5).Value2 = '  -2.82%  '

$ws.Cells.Item(38,5).Value2 = '  +0.06%  '

$ws.Cells.Item(39,4).NumberFormat = "@"
$ws.Cells.Item(39,4).Value2 = '9.051'
$ws.Cells.Item(39,4).Style = "Normal"
$ws.Cells.Item(39,5).Value2 = '  -6.66%  '

$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(40,4).Value2 = '0.06352'
$ws.Cells.Item(40,4).Style = "Normal"
$ws.Cells.Item(40,5).Value2 = '  -1.03%  '

$ws.Cells.Item(41,4).NumberFormat = "@"
$ws.Cells.Item(41,4).Value2 = '0.6460'
$ws.Cells.Item(41,4).Style = "Normal"
$ws.Cells.Item(41,5).Value2 = '  -2.21%  '

$ws.Cells.Item(42,4).NumberFormat = "@"
$ws.Cells.Item(42,4).Value2 = '11.53'
$ws.Cells.Item(42,4).Style = "Normal"
$ws.Cells.Item(42,5).Value2 = '  -2.70%  '

$ws.Cells.Item(43,4).NumberFormat = "@"
$ws.Cells.Item(43,4).Value2 = '0.1972'
$ws.Cells.Item(43,4).Style = "Normal"
$ws.Cells.Item(43,5).Value2 = '  -5.39%  '

$ws.Cells.Item(44,5).Value2 = '  +0.29%  '

$ws.Cells.Item(45,2).Value2 = 'Decentraland'
$ws.Cells.Item(45,3).Value2 = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Cells.Item(45,4).NumberFormat = "@"
$ws.Cells.Item(45,4).Value2 = '0.6163'
$ws.Cells.Item(45,4).Style = "Normal"
$ws.Cells.Item(45,5).Value2 = '  -3.57%  '

$ws.Cells.Item(46,2).Value2 = 'WEMIXTOKEN'
$ws.Cells.Item(46,3).Value2 = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(46,4).NumberFormat = "@"
$ws.Cells.Item(46,4).Value2 = '1.355'
$ws.Cells.Item(46,4).Style = "Normal"
$ws.Cells.Item(46,5).Value2 = '  +5.02%  '

$ws.Cells.Item(47,4).NumberFormat = "@"
$ws.Cells.Item(47,4).Value2 = '13.39'
$ws.Cells.Item(47,4).Style = "Normal"
$ws.Cells.Item(47,5).Value2 = '  -1.50%  '

$ws.Cells.Item(48,4).NumberFormat = "@"
$ws.Cells.Item(48,4).Value2 = '2.163'
$ws.Cells.Item(48,4).Style = "Normal"
$ws.Cells.Item(48,5).Value2 = '  -2.97%  '

$ws.Cells.Item(49,4).NumberFormat = "@"
$ws.Cells.Item(49,4).Value2 = '0.00000000352'
$ws.Cells.Item(49,4).Style = "Normal"
$ws.Cells.Item(49,5).Value2 = '  +9.78%  '

$ws.Cells.Item(50,4).NumberFormat = "@"
$ws.Cells.Item(50,4).Value2 = '3.497'
$ws.Cells.Item(50,4).Style = "Normal"
$ws.Cells.Item(50,5).Value2 = '  -1.12%  '

$ws.Cells.Item(51,4).NumberFormat = "@"
$ws.Cells.Item(51,4).Value2 = '2.176'
$ws.Cells.Item(51,4).Style = "Normal"
$ws.Cells.Item(51,5).Value2 = '  +6.81%  '
